# Fruta / hortaliza, semanal
# Inserts a new data row at row 143 (pushing existing rows 143-191 down to
# 144-192) and populates it with a new "Ají" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 143; this shifts rows 143:191 -> 144:192
# and extends the used range to A1:R192, exactly like the target diff.
$ws.Rows("143").Insert()

$ws.Cells.Item(143, 1).Value = 9
$ws.Cells.Item(143, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(143, 3).Value = "Metropolitana"
$ws.Cells.Item(143, 4).Value = 44524
$ws.Cells.Item(143, 5).Value = 13
$ws.Cells.Item(143, 6).Value = 100112021
$ws.Cells.Item(143, 7).Value = "Ají"
$ws.Cells.Item(143, 8).Value = "Americana (o)"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 16
$ws.Cells.Item(143, 11).Value = 34000
$ws.Cells.Item(143, 12).Value = 36000
$ws.Cells.Item(143, 13).Value = 35000
$ws.Cells.Item(143, 14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(143, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(143, 16).Value = 1400
$ws.Cells.Item(143, 17).Value = 25
$ws.Cells.Item(143, 18).Value = "Hortaliza"
